$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Text columns (A-D) must stay as literal text, not get auto-converted
# to dates/numbers by Excel's smart input parsing.
$ws.Cells.Item($row, 1).Value = "'2025-02-20"
$ws.Cells.Item($row, 2).Value = "'22:25:06"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "'07"

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 130102
$ws.Cells.Item($row, 6).Value = 140973
$ws.Cells.Item($row, 7).Value = 172045
$ws.Cells.Item($row, 8).Value = 154029
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146314
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192953
$ws.Cells.Item($row, 14).Value = 115267
$ws.Cells.Item($row, 15).Value = 46022
$ws.Cells.Item($row, 16).Value = 29180
$ws.Cells.Item($row, 17).Value = 67891
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48171
$ws.Cells.Item($row, 20).Value = -1

# Clear the implicit "quote prefix" style Excel attaches to cells that
# were entered with a leading apostrophe, so no style index is written
# for the new row (matching the rest of the data rows, which are
# unstyled).
$ws.Range("A80:D80").Style = "Normal"
